# Apply the "K" (strike count) value updates to column G, rows 2-26.
# These are the newly regenerated values computed from K = round(Strike#/... )
# per the commit message ("regen save_data to use K instead of Strike#").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 8
    3  = 4
    4  = 7
    5  = 5
    6  = 6
    7  = 3
    8  = 7
    9  = 5
    10 = 6
    11 = 7
    12 = 6
    13 = 3
    14 = 5
    15 = 1
    16 = 3
    17 = 1
    18 = 6
    19 = 2
    20 = 1
    21 = 7
    22 = 5
    23 = 3
    24 = 5
    25 = 3
    26 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
